# Automatizacion Completa Pinterest e Instagram
# Adds a "Status" result column (F) and per-row "Correct"/"Incorrect" flags
# (column C) to the login-check sheet, matching the success/fail state of
# each row's verification attempt.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 -> matched credentials: Correct
$ws.Range("C2").Value = "Correct"
$ws.Range("C2").Font.Name = "Calibri"
$ws.Range("C2").Font.Size = 11
$ws.Range("C2").Font.ColorIndex = 10

# Row 3 -> mismatched password: Fail
$ws.Range("F3").Value = "Fail"
$ws.Range("F3").Font.Name = "Calibri"
$ws.Range("F3").Font.Size = 11

# Row 2 -> matched credentials: Success
$ws.Range("F2").Value = "Success"
$ws.Range("F2").Font.Name = "Calibri"
$ws.Range("F2").Font.Size = 11

# Row 5 -> mismatched password: Incorrect
$ws.Range("F5").Value = "Incorrect"
$ws.Range("F5").Font.Name = "Calibri"
$ws.Range("F5").Font.Size = 11
$ws.Range("F5").Font.ColorIndex = 3

# Row 4 -> matched credentials: Success
$ws.Range("F4").Value = "Success"
$ws.Range("F4").Font.Name = "Calibri"
$ws.Range("F4").Font.Size = 11

# Row 6 -> matched credentials: Correct
$ws.Range("C6").Value = "Correct"
$ws.Range("C6").Font.Name = "Calibri"
$ws.Range("C6").Font.Size = 11
$ws.Range("C6").Font.ColorIndex = 10

# Row 6 -> mismatched password: Fail
$ws.Range("F6").Value = "Fail"
$ws.Range("F6").Font.Name = "Calibri"
$ws.Range("F6").Font.Size = 11
